$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 108.42857
$ws.Range("I6").Value = 63.166668
$ws.Range("J6").Value = 380
$ws.Range("K6").Value = 189.500004
$ws.Range("L6").Value = 1140
$ws.Range("M6").Value = -77.50000399999999
$ws.Range("N6").Value = -1364
$ws.Range("H17").Value = 4266.7656
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 4266.7656
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 12800.2968
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -13136.2968
$ws.Range("H33").Value = 666.0714
$ws.Range("I33").Value = 930.375
$ws.Range("K33").Value = 930.375
$ws.Range("M33").Value = -701.375
$ws.Range("H38").Value = 185.88889
$ws.Range("I38").Value = 21.625
$ws.Range("K38").Value = 64.875
$ws.Range("M38").Value = 307.125
$ws.Range("H40").Value = 3299
$ws.Range("I40").Value = 2271.5454
$ws.Range("J40").Value = 6124.5
$ws.Range("K40").Value = 2271.5454
$ws.Range("L40").Value = 6124.5
$ws.Range("M40").Value = -2096.5454
$ws.Range("N40").Value = -6474.5
$ws.Range("H61").Value = 64
$ws.Range("I61").Value = 64
$ws.Range("K61").Value = 192
$ws.Range("M61").Value = -20
$ws.Range("H99").Value = 22727534
$ws.Range("I99").Value = 27777978
$ws.Range("K99").Value = 83333934
$ws.Range("M99").Value = -83332436
$ws.Range("H131").Value = 4847.7666
$ws.Range("I131").Value = 1597.5883
$ws.Range("K131").Value = 4792.7649
$ws.Range("M131").Value = 247.2350999999999
$ws.Range("H132").Value = 5332.7856
$ws.Range("I132").Value = 5392.615
$ws.Range("K132").Value = 16177.845
$ws.Range("M132").Value = -13647.845
$ws.Range("H141").Value = 15288.117
$ws.Range("I141").Value = 11291.846
$ws.Range("K141").Value = 33875.538
$ws.Range("M141").Value = -28695.538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12159.189
$ws.Range("I32").Value = 7980.472
$ws.Range("J32").Value = 18997.092
$ws.Range("K32").Value = 7980.472
$ws.Range("L32").Value = 18997.092
$ws.Range("M32").Value = -7693.472
$ws.Range("N32").Value = -19571.092
$ws.Range("H34").Value = 45996
$ws.Range("I34").Value = 45996
$ws.Range("K34").Value = 45996
$ws.Range("M34").Value = -45725
$ws.Range("H54").Value = 33333.332
$ws.Range("J54").Value = 33333.332
$ws.Range("L54").Value = 33333.332
$ws.Range("N54").Value = -34871.332
$ws.Range("H61").Value = 2821.453
$ws.Range("I61").Value = 2591.8
$ws.Range("K61").Value = 2591.8
$ws.Range("M61").Value = -2379.8
$ws.Range("H136").Value = 2821.453
$ws.Range("I136").Value = 2591.8
$ws.Range("K136").Value = 7775.400000000001
$ws.Range("M136").Value = -5225.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3161.5
$ws.Range("I99").Value = 2599.9092
$ws.Range("K99").Value = 2599.9092
$ws.Range("M99").Value = -1101.9092
$ws.Range("H105").Value = 1000
$ws.Range("I105").Value = 1000
$ws.Range("K105").Value = 1000
$ws.Range("M105").Value = 747

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 866.6667
$ws.Range("H93").Value = 20674.777
$ws.Range("I93").Value = 10759.125
$ws.Range("K93").Value = 10759.125
$ws.Range("M93").Value = -8887.125
$ws.Range("H99").Value = 4170.1055
$ws.Range("I99").Value = 3841.182
$ws.Range("J99").Value = 4622.375
$ws.Range("K99").Value = 3841.182
$ws.Range("L99").Value = 4622.375
$ws.Range("M99").Value = -2343.182
$ws.Range("N99").Value = -7618.375
$ws.Range("H107").Value = 45457510
$ws.Range("I107").Value = 2277.3845
$ws.Range("J107").Value = 111115070
$ws.Range("K107").Value = 2277.3845
$ws.Range("L107").Value = 111115070
$ws.Range("M107").Value = -357.3845000000001
$ws.Range("N107").Value = -111118910
$ws.Range("H126").Value = 4170.1055
$ws.Range("I126").Value = 3841.182
$ws.Range("J126").Value = 4622.375
$ws.Range("K126").Value = 11523.546
$ws.Range("L126").Value = 13867.125
$ws.Range("M126").Value = -9053.545999999998
$ws.Range("N126").Value = -18807.125
$ws.Range("H132").Value = 53633
$ws.Range("I132").Value = 36361.863
$ws.Range("K132").Value = 109085.589
$ws.Range("M132").Value = -106555.589
$ws.Range("H134").Value = 30349.697
$ws.Range("I134").Value = 39456.168
$ws.Range("J134").Value = 6065.778
$ws.Range("K134").Value = 118368.504
$ws.Range("L134").Value = 18197.334
$ws.Range("M134").Value = -115833.504
$ws.Range("N134").Value = -23267.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H81").Value = 7616.9414
$ws.Range("I81").Value = 260
$ws.Range("J81").Value = 8597.866
$ws.Range("K81").Value = 780
$ws.Range("L81").Value = 25793.598
$ws.Range("M81").Value = 343
$ws.Range("N81").Value = -28039.598
$ws.Range("H84").Value = 7616.9414
$ws.Range("I84").Value = 260
$ws.Range("J84").Value = 8597.866
$ws.Range("K84").Value = 2340
$ws.Range("L84").Value = 77380.79399999999
$ws.Range("M84").Value = 3276
$ws.Range("N84").Value = -88612.79399999999
$ws.Range("H122").Value = 1108.1875
$ws.Range("J122").Value = 954.2857
$ws.Range("L122").Value = 8588.5713
$ws.Range("N122").Value = -13488.5713
$ws.Range("H131").Value = 11577452
$ws.Range("I131").Value = 10417568
$ws.Range("K131").Value = 31252704
$ws.Range("M131").Value = -31247664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1363.75
$ws.Range("I9").Value = 305
$ws.Range("J9").Value = 2422.5
$ws.Range("K9").Value = 305
$ws.Range("L9").Value = 2422.5
$ws.Range("M9").Value = -81
$ws.Range("N9").Value = -2870.5
$ws.Range("H30").Value = 166668640
$ws.Range("I30").Value = 200001970
$ws.Range("K30").Value = 200001970
$ws.Range("M30").Value = -200001862

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8177.756
$ws.Range("J62").Value = 8422.795
$ws.Range("L62").Value = 8422.795
$ws.Range("N62").Value = -9670.795
$ws.Range("H65").Value = 8177.756
$ws.Range("J65").Value = 8422.795
$ws.Range("L65").Value = 42113.975
$ws.Range("N65").Value = -48353.975
$ws.Range("H122").Value = 3182.75
$ws.Range("I122").Value = 2306.8462
$ws.Range("J122").Value = 4217.909
$ws.Range("K122").Value = 6920.5386
$ws.Range("L122").Value = 12653.727
$ws.Range("M122").Value = -4470.5386
$ws.Range("N122").Value = -17553.727
$ws.Range("H126").Value = 3100.889
$ws.Range("I126").Value = 3268
$ws.Range("J126").Value = 2766.6667
$ws.Range("K126").Value = 9804
$ws.Range("L126").Value = 8300.000100000001
$ws.Range("M126").Value = -7334
$ws.Range("N126").Value = -13240.0001
$ws.Range("H132").Value = 251730.48
$ws.Range("I132").Value = 4406.6
$ws.Range("K132").Value = 13219.8
$ws.Range("M132").Value = -10689.8
